$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - HANWHA AEROSPACE
$ws.Range("D2").Value = 881000
$ws.Range("E2").Value = 34.2
$ws.Range("F2").Value = 3.4
$ws.Range("H2").Value = 56
$ws.Range("K2").Value = 54.6
$ws.Range("N2").Value = 54.82400714602223

# Row 3 - HYUNDAI ROTEM
$ws.Range("D3").Value = 180700
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 2.85
$ws.Range("K3").Value = 51.8
$ws.Range("N3").Value = 54.82400714602223

# Row 4 - HANWHA SYSTEMS
$ws.Range("K4").Value = 46.4
$ws.Range("N4").Value = 54.82400714602223

# Row 5 - KOREA AEROSPACE
$ws.Range("D5").Value = 105100
$ws.Range("E5").Value = 38.4
$ws.Range("F5").Value = -3.49
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 63
$ws.Range("J5").Value = 66
$ws.Range("K5").Value = 44.6
$ws.Range("N5").Value = 54.82400714602223

# Row 6 - LIG Nex1
$ws.Range("D6").Value = 368000
$ws.Range("E6").Value = 27
$ws.Range("F6").Value = -3.92
$ws.Range("K6").Value = 37.8
$ws.Range("N6").Value = 54.82400714602223
